$d = $word.ActiveDocument

# The site rebuild dropped the trailing "Ver no Jupiter ..." / copyright
# footer block (and the blank paragraph that separated it from the
# "Requisitos" section) from the bottom of the page, while leaving the
# very last (blank, page-break) paragraph untouched.
#
# Find the paragraph that still survives ("LOQ4205: Sistemas Produtivos II
# (Requisito fraco)") and delete the three paragraphs that immediately
# follow it:
#   1) a blank "Normal" paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#      pages. Original theme under Creative Commons Attribution"

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*LOQ4205: Sistemas Produtivos II (Requisito fraco)*") {
        $target = $i
    }
}

if ($target -ne $null) {
    $startPara = $d.Paragraphs.Item($target + 1)
    $endPara = $d.Paragraphs.Item($target + 3)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
